# Add two new columns, "I0" (I) and "IF" (J), to the right of the existing
# "IP" (H) column. I0 is a constant 1 for every data row; IF mirrors the
# value already present in the corresponding H (IP) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - match the formatting already used by the other header
# cells (bold, centered, top-aligned, thin box border).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data rows 2-35: I = 1 (constant), J = same value as H (IP) on that row.
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ipValue = $ws.Cells.Item($r, 8).Formula
    $ws.Cells.Item($r, 10).Value = $ipValue
}
